$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.036.36"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.490.84"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'318.14"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").Value = "'104.83"
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").Value = "'38.73"
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").Value = "'20.03"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "2.879.72"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "2.493.32"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "'0.835"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").Value = "47.911.37"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'12.77"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("D20").Value = "'2.91"
$ws.Range("E20").Value = "  +7.38%  "
$ws.Range("D21").Value = "'6.53"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "0.0₃0928"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").Value = "'280.74"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "'70.78"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'25.64"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("E29").Value = "  -5.10%  "
$ws.Range("D30").Value = "'0.138"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "'34.46"
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").Value = "'48.99"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'19.24"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'4.49"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'119.87"
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'2.20"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "'21.73"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").Value = "1.985.89"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").Value = "'3.11"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("E47").Value = "  +6.25%  "
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "'5.09"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "'79.28"
$ws.Range("E51").Value = "  -0.51%  "
